# Auto-generated Excel COM-interop script
# Applies cell value updates described by the commit "Add data for 2023-03-29"
# to output/violent-crime-full-year.xlsx (crime statistics workbook).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 7272
$ws.Range("J2").Value = 1540
$ws.Range("J3").Value = 1616
$ws.Range("B4").Value = 1669
$ws.Range("E4").Value = 1985
$ws.Range("J4").Value = 363
$ws.Range("J5").Value = 112
$ws.Range("I6").Value = 8969
$ws.Range("J6").Value = 2108
$ws.Range("B7").Value = 23301
$ws.Range("E7").Value = 25989
$ws.Range("J7").Value = 5739

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J2").Value = 24
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J6").Value = 63
$ws.Range("J7").Value = 202

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 32
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 29
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 163
$ws.Range("J8").Value = 353
$ws.Range("J15").Value = 73
$ws.Range("J17").Value = 16
$ws.Range("J18").Value = 71
$ws.Range("J19").Value = 201
$ws.Range("J23").Value = 52
$ws.Range("J25").Value = 33
$ws.Range("J27").Value = 34
$ws.Range("J29").Value = 320
$ws.Range("J33").Value = 240
$ws.Range("J36").Value = 89
$ws.Range("J42").Value = 223
$ws.Range("J44").Value = 46
$ws.Range("J46").Value = 22
$ws.Range("J47").Value = 51
$ws.Range("J48").Value = 45
$ws.Range("J49").Value = 35
$ws.Range("J52").Value = 129
$ws.Range("J55").Value = 69
$ws.Range("J56").Value = 5
$ws.Range("J59").Value = 12
$ws.Range("B63").Value = 373
$ws.Range("E63").Value = 330
$ws.Range("J63").Value = 21
$ws.Range("J65").Value = 146
$ws.Range("J67").Value = 202
$ws.Range("J77").Value = 42
$ws.Range("J78").Value = 75
$ws.Range("J79").Value = 179
$ws.Range("J83").Value = 141
$ws.Range("J85").Value = 260
$ws.Range("J88").Value = 59
$ws.Range("J90").Value = 64
$ws.Range("J92").Value = 18
$ws.Range("J94").Value = 46
$ws.Range("J96").Value = 71
$ws.Range("B101").Value = 23301
$ws.Range("E101").Value = 25989
$ws.Range("J101").Value = 5739

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 141

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 69
$ws.Range("J4").Value = 10
$ws.Range("J6").Value = 93
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 94
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 320

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J2").Value = 50
$ws.Range("J3").Value = 56
$ws.Range("J5").Value = 13
$ws.Range("J7").Value = 201

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 63
$ws.Range("J3").Value = 104
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 43
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 18
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J3").Value = 64
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 71

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 31
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 89

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 41
$ws.Range("J7").Value = 129

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J2").Value = 8
$ws.Range("J3").Value = 7
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J2").Value = 11
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J2").Value = 5
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 18

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 16
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 118
$ws.Range("J5").Value = 9
$ws.Range("J7").Value = 353

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("I2").Value = 38
$ws.Range("I6").Value = 28

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 11
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("J3").Value = 2
$ws.Range("J7").Value = 5

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 53
$ws.Range("J6").Value = 56
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 29
